$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Propagate formatting BEFORE we touch row 220's own style ---
# Donor A: row 220 as it stands in the source file (style pattern used by rows 225 & 227).
$ws.Range("A220:G220").Copy()
$ws.Range("A225:G225").PasteSpecial(-4122)
$ws.Range("A226:G226").PasteSpecial(-4122)
$ws.Range("A227:G227").PasteSpecial(-4122)

# Donor B: row 219, the "normal" style used for the bulk of the table (and for row 220 post-edit).
$ws.Range("A219:G219").Copy()
$ws.Range("A220:G220").PasteSpecial(-4122)
$ws.Range("A221:G224").PasteSpecial(-4122)

# Overrides within rows 226/227 that differ from their row-220-style donor
# (single-cell copy/paste so the clipboard block width matches the destination).
$ws.Range("C219").Copy()
$ws.Range("C226").PasteSpecial(-4122)
$ws.Range("F219:G219").Copy()
$ws.Range("F227:G227").PasteSpecial(-4122)

# --- 2) Fill in the new daily figures (rows 221-227) ---
$ws.Range("A221").Value = 44098
$ws.Range("B221").Value = 71
$ws.Range("C221").Value = 2
$ws.Range("D221").Value = 777
$ws.Range("F221").Value = 6
$ws.Range("G221").Value = 2

$ws.Range("A222").Value = 44099
$ws.Range("B222").Value = 111
$ws.Range("C222").Value = 4
$ws.Range("D222").Value = 858
$ws.Range("F222").Value = 4
$ws.Range("G222").Value = 3

$ws.Range("A223").Value = 44100
$ws.Range("B223").Value = 82
$ws.Range("C223").Value = 3
$ws.Range("D223").Value = 851
$ws.Range("F223").Value = 8
$ws.Range("G223").Value = 4

$ws.Range("A224").Value = 44101
$ws.Range("B224").Value = 150
$ws.Range("C224").Value = 4
$ws.Range("D224").Value = 950
$ws.Range("F224").Value = 6
$ws.Range("G224").Value = 4

$ws.Range("A225").Value = 44102
$ws.Range("B225").Value = 115
$ws.Range("C225").Value = 3
$ws.Range("D225").Value = 1011
$ws.Range("F225").Value = 8
$ws.Range("G225").Value = 5

$ws.Range("A226").Value = 44103
$ws.Range("B226").Value = 101
$ws.Range("C226").Value = 4
$ws.Range("D226").Value = 1062
$ws.Range("F226").Value = 13
$ws.Range("G226").Value = 5

$ws.Range("A227").Value = 44104
$ws.Range("B227").Value = 89
$ws.Range("C227").Value = 3
$ws.Range("D227").Value = 1121
$ws.Range("F227").Value = 6
$ws.Range("G227").Value = 4

# --- 3) Extend the "active cases" formula (D-F) down through the new rows ---
$ws.Range("E221").Formula = "=D221-F221"
$ws.Range("E222").Formula = "=D222-F222"
$ws.Range("E223").Formula = "=D223-F223"
$ws.Range("E224").Formula = "=D224-F224"
$ws.Range("E225").Formula = "=D225-F225"
$ws.Range("E226").Formula = "=D226-F226"
$ws.Range("E227").Formula = "=D227-F227"
